$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): rotate C/D/E headers left (C<-D, D<-E, E<-old C)
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Update data rows 2-6: C becomes the species text (same as D), D stays the same,
# E becomes the numeric value 1
for ($r = 2; $r -le 6; $r++) {
    $species = $ws.Cells.Item($r, 4).Text
    $ws.Cells.Item($r, 3).Value = $species
    $ws.Cells.Item($r, 5).Value = 1
}
